$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(0.3895586535830836, 0.3895586535830836, 0.1926588447844201, 0.4389292024739526, 0.2145128497976859, 9)
    3  = @(0.2985262890045486, 0.3172922402961569, 0.1252830370555511, 0.3539534391068282, 0.1994532552043115, 11)
    4  = @(0.2787209351834913, 0.3278067591540666, 0.1386282432761488, 0.3723281392483635, 0.2589153759884639, 11)
    5  = @(0.3034897228953755, 0.3268934070108915, 0.1294052293206626, 0.35972938345465,   0.2064647113882407, 8)
    6  = @(0.3209127895838344, 0.3481513866787824, 0.1477021037193441, 0.3843203139561375, 0.2242915978162668, 9)
    7  = @(0.2986844067393156, 0.3331691208960183, 0.1354274215115737, 0.3680046487635362, 0.2298199212609901, 8)
    8  = @(0.3065635236640661, 0.335248462641951,  0.1402267591879591, 0.3744686357867092, 0.2322781938489615, 7)
    9  = @(0.323074747095869,  0.3439312876808158, 0.1492574881197667, 0.3863385666999435, 0.2288235751843726, 7)
    10 = @(0.3092991727900825, 0.3325527292829896, 0.1414838076039549, 0.3761433338555329, 0.2393162899682282, 5)
    11 = @(0.2921799865119131, 0.3130099903630257, 0.1357437534307621, 0.3684341914518278, 0.2591643208793996, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
